# ---------------------------------------------------------------------------
# google-saif_new.xlsx — add the "2:" (second) reference-control prefix
#
#   * urn_prefix_content gains a new row describing prefix "2" as the
#     reference_controls base urn.
#   * assessment_content's "reference_controls" column (F) is re-pointed
#     from the "1:GSC-xx" (threats) prefix to the new "2:GSC-xx" prefix.
#   * assorted view-state (zoom / selection / active tab) left behind by
#     the author while making the edit in the desktop app.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) urn_prefix_content: register prefix "2" -> reference_control base urn
# ---------------------------------------------------------------------------
$wsPrefix = $wb.Worksheets.Item("urn_prefix_content")
$wsPrefix.Cells.Item(3, 1).Value = 2
$wsPrefix.Cells.Item(3, 2).Value = "urn:intuitem:risk:reference_control:google-saif"
$wsPrefix.Columns.Item(2).ColumnWidth = 36.66666666666667

# ---------------------------------------------------------------------------
# 2) assessment_content: reference_controls column (F) now uses prefix 2
#    instead of prefix 1 for every GSC-xx reference control code.
# ---------------------------------------------------------------------------
$wsAssessment = $wb.Worksheets.Item("assessment_content")
for ($r = 2; $r -le 13; $r++) {
    $cell = $wsAssessment.Cells.Item($r, 6)
    $old = $cell.Text
    $new = $old.Replace("1:GSC", "2:GSC")
    $cell.Value = $new
}

# ---------------------------------------------------------------------------
# 3) Leftover view state (zoom/selection/active sheet) from the editing
#    session. Order matters: the last sheet activated ends up as the
#    workbook's active tab, matching assessment_content being active last.
# ---------------------------------------------------------------------------

# reference_controls_meta
$wsRefMeta = $wb.Worksheets.Item("reference_controls_meta")
$wsRefMeta.Activate()
$excel.ActiveWindow.Zoom = 161

# reference_controls_content
$wsRefContent = $wb.Worksheets.Item("reference_controls_content")
$wsRefContent.Activate()
$wsRefContent.Range("C33").Select()
$excel.ActiveWindow.Zoom = 163

# threats_meta
$wsThreatsMeta = $wb.Worksheets.Item("threats_meta")
$wsThreatsMeta.Activate()
$wsThreatsMeta.Range("B2").Select()
$excel.ActiveWindow.Zoom = 186

# urn_prefix_content
$wsPrefix.Activate()
$wsPrefix.Range("B21").Select()
$excel.ActiveWindow.Zoom = 185

# assessment_content (activated last -> becomes the active tab)
$wsAssessment.Activate()
$wsAssessment.Range("F1:F1048576").Select()
$excel.ActiveWindow.Zoom = 156

"done"
